$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "fecha" column (B) holds plain text dates like "2026-02-13" in the
# existing rows (stored as text, not real Excel dates). Pre-format the
# new range as Text so Excel doesn't auto-convert the yyyy-mm-dd strings
# into date serials when they're assigned below.
$ws.Range("B21:B31").NumberFormat = "@"

$responsable = '{"nombre":"Luisa","rol":"ADMIN"}'

$rows = @(
    @{ r=21; id=1771624904383; fecha="2026-02-20"; desc="Cafeteria Doña Martha"; cat="Insumos";   monto=632240 },
    @{ r=22; id=1771624963201; fecha="2026-02-20"; desc="Agua";                    cat="Servicios"; monto=432850 },
    @{ r=23; id=1771632832729; fecha="2026-02-20"; desc="Compras";                 cat="Insumos";   monto=2224290 },
    @{ r=24; id=1771633930022; fecha="2026-02-20"; desc="Zumo de limon";           cat="Insumos";   monto=40000 },
    @{ r=25; id=1771633942158; fecha="2026-02-20"; desc="Martin";                  cat="Nomina";    monto=50000 },
    @{ r=26; id=1771639250275; fecha="2026-02-20"; desc="Bailarin";                cat="Nomina";    monto=20000 },
    @{ r=27; id=1771639250285; fecha="2026-02-20"; desc="Bailarin";                cat="Nomina";    monto=20000 },
    @{ r=28; id=1771641168472; fecha="2026-02-20"; desc="Hielo";                   cat="Insumos";   monto=84000 },
    @{ r=29; id=1771644118288; fecha="2026-02-20"; desc="Luz Helena";              cat="Nomina";    monto=80000 },
    @{ r=30; id=1771644127448; fecha="2026-02-20"; desc="Fredy";                   cat="Nomina";    monto=80000 },
    @{ r=31; id=1771644141358; fecha="2026-02-20"; desc="Luisa ";                  cat="Nomina";    monto=80000 }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.id
    $ws.Cells.Item($r, 2).Value = $row.fecha
    $ws.Cells.Item($r, 3).Value = $row.desc
    $ws.Cells.Item($r, 4).Value = $row.cat
    $ws.Cells.Item($r, 5).Value = $row.monto
    $ws.Cells.Item($r, 6).Value = $responsable
}
